# Updates to tech selection
# Applies data + view changes to the "Updates" worksheet of the techs_database workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Updates")
$ws.Activate()

# --- Cell value updates -----------------------------------------------------
# Column N ("Cycle life (# of cycles)") updates, per-row.
$nValues = @{
    2  = 5000
    3  = 5000
    4  = 5000
    5  = 5000
    6  = 5000
    7  = 2000
    8  = 2000
    9  = 2000
    10 = 2000
    11 = 2000
    12 = 3000
    13 = 3000
    14 = 3000
    15 = 3000
    16 = 3000
    17 = 3000
    18 = 3000
    19 = 3000
    20 = 3000
}

foreach ($row in $nValues.Keys) {
    $ws.Range("N$row").Value2 = $nValues[$row]
}

# Rows 12 and 13 also get updated feasibility scores in columns W and X.
$ws.Range("W12").Value2 = 1
$ws.Range("X12").Value2 = 1
$ws.Range("W13").Value2 = 1
$ws.Range("X13").Value2 = 1

# Columns Y and Z ("Feas. score for distribution" / "Feas. score for
# transmission") become 1 for rows 2-20, matching the style already used by
# columns W/X in the same row (general alignment instead of right-aligned).
for ($row = 2; $row -le 20; $row++) {
    $yCell = $ws.Range("Y$row")
    $zCell = $ws.Range("Z$row")
    $yCell.Value2 = 1
    $zCell.Value2 = 1
    $yCell.Style = $ws.Range("W$row").Style
    $zCell.Style = $ws.Range("W$row").Style
}

# --- New data point -----------------------------------------------------
# Flywheel - Short duration (row 50) gains a discharge duration value.
$ws.Range("E50").Value2 = 0.5

# --- View / selection updates -----------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("E51").Select() | Out-Null
